# Auto-generated edit script applying the diff to Excalibur_Profits workbook
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 6464.2
$ws.Range("I18").Value = 1129.4
$ws.Range("K18").Value = 1129.4
$ws.Range("M18").Value = -845.4000000000001

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H132").Value = 51364.13
$ws.Range("I132").Value = 52437.26
$ws.Range("K132").Value = 157311.78
$ws.Range("M132").Value = -154781.78

$ws.Range("H135").Value = 1504.4333
$ws.Range("I135").Value = 1540.5
$ws.Range("J135").Value = 999.5
$ws.Range("K135").Value = 13864.5
$ws.Range("L135").Value = 8995.5
$ws.Range("M135").Value = -11329.5
$ws.Range("N135").Value = -14065.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 6347.25
$ws.Range("I63").Value = 2694.5
$ws.Range("K63").Value = 2694.5
$ws.Range("M63").Value = -2008.5

$ws.Range("H66").Value = 6347.25
$ws.Range("I66").Value = 2694.5
$ws.Range("K66").Value = 13472.5
$ws.Range("M66").Value = -10040.5

$ws.Range("H86").Value = 90251.2
$ws.Range("J86").Value = 90251.2
$ws.Range("L86").Value = 90251.2
$ws.Range("N86").Value = -92623.2

$ws.Range("H89").Value = 90251.2
$ws.Range("J89").Value = 90251.2
$ws.Range("L89").Value = 270753.6
$ws.Range("N89").Value = -282609.6

$ws.Range("H95").Value = 32019.6
$ws.Range("J95").Value = 32019.6
$ws.Range("L95").Value = 32019.6
$ws.Range("N95").Value = -37511.6

$ws.Range("H132").Value = 667967.1
$ws.Range("I132").Value = 804063.1
$ws.Range("K132").Value = 2412189.3
$ws.Range("M132").Value = -2409659.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1863.6818
$ws.Range("I86").Value = 1798.25
$ws.Range("K86").Value = 1798.25
$ws.Range("M86").Value = -675.25

$ws.Range("H88").Value = 38246.816
$ws.Range("J88").Value = 38246.816
$ws.Range("L88").Value = 38246.816
$ws.Range("N88").Value = -39058.816

$ws.Range("H89").Value = 1863.6818
$ws.Range("I89").Value = 1798.25
$ws.Range("K89").Value = 8991.25
$ws.Range("M89").Value = -3375.25

$ws.Range("H91").Value = 38246.816
$ws.Range("J91").Value = 38246.816
$ws.Range("L91").Value = 38246.816
$ws.Range("N91").Value = -41054.816

$ws.Range("H134").Value = 580914
$ws.Range("I134").Value = 777127.4
$ws.Range("J134").Value = 10111.546
$ws.Range("K134").Value = 2331382.2
$ws.Range("L134").Value = 30334.638
$ws.Range("M134").Value = -2328847.2
$ws.Range("N134").Value = -35404.638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H22").Value = 742.5789
$ws.Range("I22").Value = 700.3125
$ws.Range("J22").Value = 968
$ws.Range("K22").Value = 700.3125
$ws.Range("L22").Value = 968
$ws.Range("M22").Value = -350.3125
$ws.Range("N22").Value = -1668

$ws.Range("H132").Value = 22761206
$ws.Range("I132").Value = 45407.125
$ws.Range("J132").Value = 83336664
$ws.Range("K132").Value = 136221.375
$ws.Range("L132").Value = 250009992
$ws.Range("M132").Value = -133691.375
$ws.Range("N132").Value = -250015052

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 126692200
$ws.Range("I4").Value = 250133920
$ws.Range("K4").Value = 750401760
$ws.Range("M4").Value = -750401648

$ws.Range("H25").Value = 25
$ws.Range("I25").Value = 25
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 75
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 94
$ws.Range("N25").ClearContents()

$ws.Range("H30").Value = 25
$ws.Range("I30").Value = 25
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 75
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 27
$ws.Range("N30").ClearContents()

$ws.Range("H31").Value = 666.6667
$ws.Range("I31").Value = 666.6667
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2000.0001
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1712.0001
$ws.Range("N31").ClearContents()

$ws.Range("H131").Value = 9660.186
$ws.Range("I131").Value = 642.7
$ws.Range("J131").Value = 14964.588
$ws.Range("K131").Value = 1928.1
$ws.Range("L131").Value = 44893.764
$ws.Range("M131").Value = 3111.9
$ws.Range("N131").Value = -54973.764

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()

$ws.Range("H97").Value = 2395.1943
$ws.Range("I97").Value = 794.5789
$ws.Range("J97").Value = 4184.1177
$ws.Range("K97").Value = 794.5789
$ws.Range("L97").Value = 4184.1177
$ws.Range("M97").Value = -298.5789
$ws.Range("N97").Value = -5176.1177

$ws.Range("H132").Value = 755062.0600000001
$ws.Range("I132").Value = 862406.3
$ws.Range("K132").Value = 2587218.9
$ws.Range("M132").Value = -2584688.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 900.43475
$ws.Range("J46").Value = 926.6667
$ws.Range("L46").Value = 926.6667
$ws.Range("N46").Value = -1302.6667

$ws.Range("H58").Value = 9154.75
$ws.Range("I58").Value = 8046.5
$ws.Range("K58").Value = 8046.5
$ws.Range("M58").Value = -7786.5

$ws.Range("H64").Value = 28285.4
$ws.Range("I64").Value = 1431
$ws.Range("J64").Value = 34999
$ws.Range("K64").Value = 1431
$ws.Range("L64").Value = 34999
$ws.Range("M64").Value = -1206
$ws.Range("N64").Value = -35449

$ws.Range("H67").Value = 28285.4
$ws.Range("I67").Value = 1431
$ws.Range("J67").Value = 34999
$ws.Range("K67").Value = 1431
$ws.Range("L67").Value = 34999
$ws.Range("M67").Value = -651
$ws.Range("N67").Value = -36559

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0

$ws.Range("H132").Value = 1155210.4
$ws.Range("I132").Value = 1194838.2
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 3584514.6
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -3581984.6
$ws.Range("N132").Value = -23060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

$ws.Range("H63").Value = 38061.93
$ws.Range("I63").Value = 18000
$ws.Range("J63").Value = 39605.152
$ws.Range("K63").Value = 18000
$ws.Range("L63").Value = 39605.152
$ws.Range("M63").Value = -17376
$ws.Range("N63").Value = -40853.152

$ws.Range("H66").Value = 38061.93
$ws.Range("I66").Value = 18000
$ws.Range("J66").Value = 39605.152
$ws.Range("K66").Value = 54000
$ws.Range("L66").Value = 118815.456
$ws.Range("M66").Value = -50880
$ws.Range("N66").Value = -125055.456

$ws.Range("H75").Value = 79722
$ws.Range("J75").Value = 79722
$ws.Range("L75").Value = 79722
$ws.Range("N75").Value = -81594

$ws.Range("H78").Value = 79722
$ws.Range("J78").Value = 79722
$ws.Range("L78").Value = 239166
$ws.Range("N78").Value = -248526

$ws.Range("H132").Value = 6101188
$ws.Range("I132").Value = 6710740
$ws.Range("J132").Value = 5665
$ws.Range("K132").Value = 20132220
$ws.Range("L132").Value = 16995
$ws.Range("M132").Value = -20129690
$ws.Range("N132").Value = -22055
